$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, shifting existing rows 32-36 down to 33-37
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly record
$ws.Cells.Item(32, 1).Value = 9
$ws.Cells.Item(32, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(32, 3).Value = "Metropolitana"
$ws.Cells.Item(32, 4).Value = 44452
$ws.Cells.Item(32, 5).Value = 13
$ws.Cells.Item(32, 6).Value = "Fruta"
$ws.Cells.Item(32, 7).Value = 100108
$ws.Cells.Item(32, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(32, 9).Value = 100108003
$ws.Cells.Item(32, 10).Value = "Maracuyá"
$ws.Cells.Item(32, 11).Value = "Sin especificar"
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 20
$ws.Cells.Item(32, 14).Value = 36000
$ws.Cells.Item(32, 15).Value = 36000
$ws.Cells.Item(32, 16).Value = 36000
$ws.Cells.Item(32, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(32, 18).Value = "Perú"
$ws.Cells.Item(32, 19).Value = 2000
$ws.Cells.Item(32, 20).Value = 18
